$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New attachment rows scraped since the last update (rows 76-82) ---

# Row 76
$ws.Range("A76").Value = "2025-04-18 10:26"
$ws.Range("B76").Value = "http://www.scnj.gov.cn/oldfiles/njxxxgk/2016/12/04/20161204095941-621511.xls"
$ws.Range("C76").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=dca978027c85406daf5679fa59b345e8&type=0"
$ws.Range("D76").Value = "老干局"
$ws.Range("E76").Value = "老干部局"
$ws.Range("F76").Value = "附件：2015年老干局决算批复表.xls"
$ws.Range("G76").Value = "http://www.scnj.gov.cn/public/6598311/11926851.html"
$ws.Range("H76").Value = "中共南江县委老干部局2015年部门决算"

# Row 77
$ws.Range("A77").Value = "2025-04-18 10:26"
$ws.Range("B77").Value = "http://www.scnj.gov.cn/oldfiles/njxxxgk/2016/02/29/20160229205551-832962.pdf"
$ws.Range("C77").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=0727c0d393174ab8858320f4f2426b5d&type=0"
$ws.Range("D77").Value = "卯足干劲"
$ws.Range("E77").Value = "铆足干劲"
$ws.Range("F77").Value = "刘凯同志在全县招商引资暨回乡创业人士座谈会上的讲话"
$ws.Range("G77").Value = "http://www.scnj.gov.cn/public/6598711/12275741.html"
$ws.Range("H77").Value = "刘凯同志在全县招商引资暨回乡创业人士座谈会上的讲话"

# Row 78
$ws.Range("A78").Value = "2025-04-21 08:36"
$ws.Range("B78").Value = "http://www.scnj.gov.cn/group3/M00/06/A5/rBUtImF4162ASr8rAEMIidR7yHs615.pdf"
$ws.Range("C78").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=acc82975adec403e8b5813263d7915a9&type=0"
$ws.Range("D78").Value = "党的十九大及十九届二中、三中、四中、五中全会精神"
$ws.Range("E78").Value = "党的十九大和十九届二中、三中、四中、五中全会精神"
$ws.Range("F78").Value = "附件.pdf"
$ws.Range("G78").Value = "http://www.scnj.gov.cn/public/6598011/13433071.html"
$ws.Range("H78").Value = "巴中市人民政府关于印发《巴中市“十四五”水安全保障规划》的通知"

# Row 79
$ws.Range("A79").Value = "2025-04-21 08:36"
$ws.Range("B79").Value = "http://www.scnj.gov.cn/oldfiles/njxxxgk/2019/12/19/20191219105746-336359.pdf"
$ws.Range("C79").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=87ad10a7f25f4a13a354ebebd455ef22&type=0"
$ws.Range("D79").Value = "建档立卡贫因户"
$ws.Range("E79").Value = "建档立卡贫困户"
$ws.Range("F79").Value = "关于进一步事实无人抚养儿童保障工作的通知"
$ws.Range("G79").Value = "http://www.scnj.gov.cn/public/6599051/12433631.html"
$ws.Range("H79").Value = "关于进一步事实无人抚养儿童保障工作的通知"

# Row 80
$ws.Range("A80").Value = "2025-04-22 09:19"
$ws.Range("B80").Value = "http://www.scnj.gov.cn/oldfiles/njx/file/p/f8695f791e669c99caf0f47879381f72.doc"
$ws.Range("C80").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=9162d69f04634592b61fd1e34a5b7c57&type=0"
$ws.Range("D80").Value = "中华人民共和国行政处罚法》"
$ws.Range("E80").Value = "《中华人民共和国行政处罚法》"
$ws.Range("F80").Value = "f8695f791e669c99caf0f47879381f72.doc"
$ws.Range("G80").Value = "http://www.scnj.gov.cn/xxgk/wgk/glgk/12587231.html"

# Row 81
$ws.Range("A81").Value = "2025-04-24 10:37"
$ws.Range("B81").Value = "http://www.scnj.gov.cn/group3/M00/04/4F/rBUtImB1Cg6AMRI5ADpQBNBWy24062.pdf"
$ws.Range("C81").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=318060029b594ab9996060ac8edae011&type=0"
$ws.Range("D81").Value = "党的十九大及十九届二中、三中、四中、五中全会精神"
$ws.Range("E81").Value = "党的十九大和十九届二中、三中、四中、五中全会精神"
$ws.Range("F81").Value = "关于同意《县应急管理局2021年安全生产监管监察执法工作计划》的批复.pdf"
$ws.Range("G81").Value = "http://www.scnj.gov.cn/public/6598411/13148181.html"
$ws.Range("H81").Value = "南江县应急管理局2021年安全生产监管监察执法工作计划"

# Row 82
$ws.Range("A82").Value = "2025-04-24 15:05"
$ws.Range("B82").Value = "http://www.scnj.gov.cn/group3/M00/0B/8C/rBUtImN1uK2ATxgQAAB-ADb6Cog952.doc"
$ws.Range("C82").Value = "http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=fee08a86788e43e4b59be5cc5029f93f&type=0"
$ws.Range("D82").Value = "辩认"
$ws.Range("E82").Value = "辨认"
$ws.Range("F82").Value = "四川省住房和城乡建设行政处罚裁量权适用规定.doc"
$ws.Range("G82").Value = "http://www.scnj.gov.cn/public/6598631/13804095.html"
$ws.Range("H82").Value = "四川省住房和城乡建设厅关于印发《四川省住房和城乡建设行政处罚裁量权适用规定》《四川省住房和城乡建设行政处罚裁量标准》的通知"

# Match the workbook's existing body-row formatting (vertical-centered) for
# row 75 (previously the last row, now re-touched) through the new rows
# 76-82, consistent with every other data row in the sheet.
$ws.Range("A75:H82").VerticalAlignment = -4108
